$d = $word.ActiveDocument

# Start from the current last paragraph ("Explicando a diferença entre Git e GitHub")
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Last
$newPara1.Range.Text = "Prática"

$newPara2Anchor = $d.Paragraphs.Last.Range
$newPara2Anchor.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Last
$newPara2.Range.Text = "Introdução ao GitFlow"
